$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force numeric-looking price strings to stay text (matches source inlineStr type)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '97.246.41'
$ws.Range("E2").Value = '  +0.31%  '
$ws.Range("D3").Value = '3.713.00'
$ws.Range("E3").Value = '  +1.18%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '1.99'
$ws.Range("E5").Value = '  +4.74%  '
$ws.Range("D6").Value = '236.45'
$ws.Range("E6").Value = '  -1.22%  '
$ws.Range("D7").Value = '657.23'
$ws.Range("E7").Value = '  +0.48%  '
$ws.Range("E8").Value = '  +2.63%  '
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").Value = '1.07'
$ws.Range("E9").Value = '  -1.08%  '
$ws.Range("B10").Value = 'USDC'
$ws.Range("C10").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D10").Value = '0.999'
$ws.Range("E10").Value = '  -0.05%  '
$ws.Range("D11").Value = '3.713.09'
$ws.Range("E11").Value = '  +1.21%  '
$ws.Range("E12").Value = '  +16.60%  '
$ws.Range("D13").Value = '44.78'
$ws.Range("E13").Value = '  -1.53%  '
$ws.Range("E14").Value = '  +0.34%  '
$ws.Range("D15").Value = '6.91'
$ws.Range("E15").Value = '  +0.43%  '
$ws.Range("D16").Value = '4.408.36'
$ws.Range("E16").Value = '  +1.24%  '
$ws.Range("D17").Value = '97.024.75'
$ws.Range("E17").Value = '  +0.39%  '
$ws.Range("D18").Value = '9.10'
$ws.Range("E18").Value = '  +0.85%  '
$ws.Range("D19").Value = '3.716.22'
$ws.Range("E19").Value = '  +1.93%  '
$ws.Range("D20").Value = '13.02'
$ws.Range("E20").Value = '  +1.96%  '
$ws.Range("D21").Value = '18.77'
$ws.Range("E21").Value = '  -1.16%  '
$ws.Range("E22").Value = '  -3.32%  '
$ws.Range("D23").Value = '525.76'
$ws.Range("E23").Value = '  -1.08%  '
$ws.Range("D24").Value = '3.47'
$ws.Range("E24").Value = '  -0.96%  '
$ws.Range("D25").Value = '0.0000226'
$ws.Range("E25").Value = '  +10.87%  '
$ws.Range("E26").Value = '  -3.78%  '
$ws.Range("D27").Value = '107.07'
$ws.Range("E27").Value = '  +4.41%  '
$ws.Range("D28").Value = '0.195'
$ws.Range("E28").Value = '  +15.62%  '
$ws.Range("D29").Value = '3.917.13'
$ws.Range("E29").Value = '  +1.27%  '
$ws.Range("D30").Value = '13.52'
$ws.Range("E30").Value = '  +0.24%  '
$ws.Range("D31").Value = '12.63'
$ws.Range("E31").Value = '  +0.91%  '
$ws.Range("D32").Value = '3.04'
$ws.Range("E32").Value = '  -0.32%  '
$ws.Range("E33").Value = '  -0.10%  '
$ws.Range("E34").Value = '  +3.21%  '
$ws.Range("E35").Value = '  -3.13%  '
$ws.Range("B36").Value = 'EthereumClassic'
$ws.Range("C36").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D36").Value = '32.73'
$ws.Range("E36").Value = '  +0.07%  '
$ws.Range("B37").Value = 'Binance-PegBSC-USD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D37").Value = '0.998'
$ws.Range("E37").Value = '  -0.17%  '
$ws.Range("D38").Value = '642.93'
$ws.Range("E38").Value = '  -2.09%  '
$ws.Range("D39").Value = '0.595'
$ws.Range("E39").Value = '  -2.00%  '
$ws.Range("D40").Value = '8.75'
$ws.Range("E40").Value = '  -2.43%  '
$ws.Range("D42").Value = '0.507'
$ws.Range("E42").Value = '  +12.24%  '
$ws.Range("D43").Value = '0.166'
$ws.Range("E43").Value = '  +2.58%  '
$ws.Range("E44").Value = '  -2.46%  '
$ws.Range("D45").Value = '2.03'
$ws.Range("E45").Value = '  +1.70%  '
$ws.Range("D46").Value = '40.35'
$ws.Range("E46").Value = '  +4.90%  '
$ws.Range("D47").Value = '0.965'
$ws.Range("E47").Value = '  +0.21%  '
$ws.Range("D48").Value = '0.0459'
$ws.Range("E48").Value = '  -0.67%  '
$ws.Range("D49").Value = '2.41'
$ws.Range("E49").Value = '  +3.32%  '
$ws.Range("D51").Value = '8.70'
$ws.Range("E51").Value = '  -0.77%  '
